$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E11: ДЗ_3 grade corrected from 0 to 5 (now "correct", so it loses the
# green highlight fill -> matches style already used by G11/H11 in this row)
$ws.Range("G11").Copy()
$ws.Range("E11").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E11").Value = 5

# Row 13: all four homework grades corrected from 0 to 5, same style change
$ws.Range("G13").Copy()
$ws.Range("C13:F13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 5

# Move the active selection to reflect the last edited cell
[void]$ws.Range("E11").Select()
